# Backlog housing.xlsx - add new "propriétaire" (owner) backlog items to Sprint 3,
# downgrade priority of two existing Sprint-3 rows, and widen a couple of columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Downgrade priority of two existing Sprint 3 rows (Haute -> Moyenne) ---
$ws.Range("C11").Value = "Moyenne"
$ws.Range("C12").Value = "Moyenne"

# --- 2. Insert 4 new backlog rows (13-16), copying formatting from an existing
#        data row so the same fills/borders/fonts get reused. ---

# Rows 13-15 look exactly like the existing wrap-text data rows.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E15").PasteSpecial(-4122)

$ws.Range("A13").Value = "Sprint 3"
$ws.Range("B13").Value = "Ajouter un logement (propriétaire)"
$ws.Range("C13").Value = "Haute"
$ws.Range("D13").Value = "À faire"
$ws.Range("E13").Value = "Formulaire avec nom, adresse, photos, prix, équipements"

$ws.Range("A14").Value = "Sprint 3"
$ws.Range("B14").Value = "Modifier un logement (propriétaire)"
$ws.Range("C14").Value = "Haute"
$ws.Range("D14").Value = "À faire"
$ws.Range("E14").Value = "Modifier informations existantes"

$ws.Range("A15").Value = "Sprint 3"
$ws.Range("B15").Value = "Supprimer un logement (propriétaire)"
$ws.Range("C15").Value = "Haute"
$ws.Range("D15").Value = "À faire"
$ws.Range("E15").Value = "Supprimer logement de la liste"

# Row 16 is the same palette but without wrap text, and column B gets a plain
# light-gray fill instead of the usual green tint.
$ws.Range("A2:E2").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A16:E16").WrapText = $False
$ws.Range("B16").Interior.ThemeColor = 2

$ws.Range("A16").Value = "Sprint 3"
$ws.Range("B16").Value = "Consulter les réservations de mes logements (propriétaire)"
$ws.Range("C16").Value = "Haute"
$ws.Range("D16").Value = "À faire"
$ws.Range("E16").Value = "Afficher réservations avec date, utilisateur, logement réservé"

$ws.Application.CutCopyMode = $False

# --- 3. Widen columns B and E to fit the new, longer text. ---
$ws.Columns.Item(2).ColumnWidth = 49.3
$ws.Columns.Item(5).ColumnWidth = 54.0

# --- 4. Autofit every data row's height (diff drops all explicit row heights). ---
$ws.Range("A1:E16").EntireRow.AutoFit()

# --- 5. Update selection / scroll position to match the post-edit view. ---
$ws.Range("E20").Select()
